$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 1.225746
$ws.Range("H2").Value = 2.451492
$ws.Range("I2").Value = 0.1045666790027646
$ws.Range("J2").Value = 0.09912108042370671
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 64.20135099999999
$ws.Range("N2").Value = 128.402702
$ws.Range("O2").Value = 0.4070144914449589
$ws.Range("P2").Value = 0.3181813759721767
$ws.Range("Q2").Value = 78.69454918284599
$ws.Range("R2").Value = 314.778196731384
$ws.Range("S2").Value = 0.04256015367639851
$ws.Range("T2").Value = 0.03153848175706379

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 1.225746
$ws.Range("H3").Value = 2.451492
$ws.Range("I3").Value = 0.1045666790027646
$ws.Range("J3").Value = 0.09912108042370671
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 15.75734066666666
$ws.Range("N3").Value = 47.27202199999999
$ws.Range("O3").Value = 0.09989612209201491
$ws.Range("P3").Value = 0.11713987922892
$ws.Range("Q3").Value = 19.314497292804
$ws.Range("R3").Value = 115.886983756824
$ws.Range("S3").Value = 0.01044580573241671
$ws.Range("T3").Value = 0.01161103138987307

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 1.225746
$ws.Range("H4").Value = 2.451492
$ws.Range("I4").Value = 0.1045666790027646
$ws.Range("J4").Value = 0.09912108042370671
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 24.57775566666666
$ws.Range("N4").Value = 73.733267
$ws.Range("O4").Value = 0.1558145205313015
$ws.Range("P4").Value = 0.1827107372630203
$ws.Range("Q4").Value = 30.126085697394
$ws.Range("R4").Value = 180.756514184364
$ws.Range("S4").Value = 0.01629300695236628
$ws.Range("T4").Value = 0.01811048568252258

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 1.225746
$ws.Range("H5").Value = 2.451492
$ws.Range("I5").Value = 0.1045666790027646
$ws.Range("J5").Value = 0.09912108042370671
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 14.32600733333333
$ws.Range("N5").Value = 42.978022
$ws.Range("O5").Value = 0.09082196088386706
$ws.Range("P5").Value = 0.1064993646046676
$ws.Range("Q5").Value = 17.560046184804
$ws.Range("R5").Value = 105.360277108824
$ws.Range("S5").Value = 0.00949695083014497
$ws.Range("T5").Value = 0.01055633208405292

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 1.225746
$ws.Range("H6").Value = 2.451492
$ws.Range("I6").Value = 0.1045666790027646
$ws.Range("J6").Value = 0.09912108042370671
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 33.41628466666666
$ws.Range("N6").Value = 100.248854
$ws.Range("O6").Value = 0.2118477555025799
$ws.Range("P6").Value = 0.2484162545532246
$ws.Range("Q6").Value = 40.959877265028
$ws.Range("R6").Value = 245.759263590168
$ws.Range("S6").Value = 0.02215221624709443
$ws.Range("T6").Value = 0.02462328754612618

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("G7").Value = 1.225746
$ws.Range("H7").Value = 2.451492
$ws.Range("I7").Value = 0.1045666790027646
$ws.Range("J7").Value = 0.09912108042370671
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 5.4585215
$ws.Range("N7").Value = 10.917043
$ws.Range("O7").Value = 0.03460514954527787
$ws.Range("P7").Value = 0.02705238837799083
$ws.Range("Q7").Value = 6.690760894538999
$ws.Range("R7").Value = 26.763043578156
$ws.Range("S7").Value = 0.003618545564343737
$ws.Range("T7").Value = 0.002681461964068177

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.6169616666666666
$ws.Range("H8").Value = 1.850885
$ws.Range("I8").Value = 0.05263213794321498
$ws.Range("J8").Value = 0.07483676101738548
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 64.20135099999999
$ws.Range("N8").Value = 128.402702
$ws.Range("O8").Value = 0.4070144914449589
$ws.Range("P8").Value = 0.3181813759721767
$ws.Range("Q8").Value = 39.60977251521166
$ws.Range("R8").Value = 237.6586350912699
$ws.Range("S8").Value = 0.02142204285861857
$ws.Range("T8").Value = 0.02381166359381267

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.6169616666666666
$ws.Range("H9").Value = 1.850885
$ws.Range("I9").Value = 0.05263213794321498
$ws.Range("J9").Value = 0.07483676101738548
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 15.75734066666666
$ws.Range("N9").Value = 47.27202199999999
$ws.Range("O9").Value = 0.09989612209201491
$ws.Range("P9").Value = 0.11713987922892
$ws.Range("Q9").Value = 9.721675159941109
$ws.Range("R9").Value = 87.49507643946998
$ws.Range("S9").Value = 0.005257746477939174
$ws.Range("T9").Value = 0.008766369147460085

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.6169616666666666
$ws.Range("H10").Value = 1.850885
$ws.Range("I10").Value = 0.05263213794321498
$ws.Range("J10").Value = 0.07483676101738548
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 24.57775566666666
$ws.Range("N10").Value = 73.733267
$ws.Range("O10").Value = 0.1558145205313015
$ws.Range("P10").Value = 0.1827107372630203
$ws.Range("Q10").Value = 15.16353309903278
$ws.Range("R10").Value = 136.471797891295
$ws.Range("S10").Value = 0.008200851338159362
$ws.Range("T10").Value = 0.01367347977986296

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.6169616666666666
$ws.Range("H11").Value = 1.850885
$ws.Range("I11").Value = 0.05263213794321498
$ws.Range("J11").Value = 0.07483676101738548
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 14.32600733333333
$ws.Range("N11").Value = 42.978022
$ws.Range("O11").Value = 0.09082196088386706
$ws.Range("P11").Value = 0.1064993646046676
$ws.Range("Q11").Value = 8.838597361052221
$ws.Range("R11").Value = 79.54737624946999
$ws.Range("S11").Value = 0.004780153973512966
$ws.Range("T11").Value = 0.007970067497422911

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.6169616666666666
$ws.Range("H12").Value = 1.850885
$ws.Range("I12").Value = 0.05263213794321498
$ws.Range("J12").Value = 0.07483676101738548
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 33.41628466666666
$ws.Range("N12").Value = 100.248854
$ws.Range("O12").Value = 0.2118477555025799
$ws.Range("P12").Value = 0.2484162545532246
$ws.Range("Q12").Value = 20.61656668175444
$ws.Range("R12").Value = 185.54910013579
$ws.Range("S12").Value = 0.01115000029057227
$ws.Range("T12").Value = 0.01859066787483367

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.6169616666666666
$ws.Range("H13").Value = 1.850885
$ws.Range("I13").Value = 0.05263213794321498
$ws.Range("J13").Value = 0.07483676101738548
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 5.4585215
$ws.Range("N13").Value = 10.917043
$ws.Range("O13").Value = 0.03460514954527787
$ws.Range("P13").Value = 0.02705238837799083
$ws.Range("Q13").Value = 3.367698522175833
$ws.Range("R13").Value = 20.206191133055
$ws.Range("S13").Value = 0.001821343004412648
$ws.Range("T13").Value = 0.002024513123993196

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.6710410000000001
$ws.Range("H14").Value = 2.013123
$ws.Range("I14").Value = 0.0572455703259029
$ws.Range("J14").Value = 0.08139652374383181
$ws.Range("K14").Value = 2
$ws.Range("M14").Value = 64.20135099999999
$ws.Range("N14").Value = 128.402702
$ws.Range("O14").Value = 0.4070144914449589
$ws.Range("P14").Value = 0.3181813759721767
$ws.Range("Q14").Value = 43.081738776391
$ws.Range("R14").Value = 258.490432658346
$ws.Range("S14").Value = 0.023299776693674
$ws.Range("T14").Value = 0.02589885792416435

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.6710410000000001
$ws.Range("H15").Value = 2.013123
$ws.Range("I15").Value = 0.0572455703259029
$ws.Range("J15").Value = 0.08139652374383181
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 15.75734066666666
$ws.Range("N15").Value = 47.27202199999999
$ws.Range("O15").Value = 0.09989612209201491
$ws.Range("P15").Value = 0.11713987922892
$ws.Range("Q15").Value = 10.57382163830067
$ws.Range("R15").Value = 95.164394744706
$ws.Range("S15").Value = 0.005718610482503422
$ws.Range("T15").Value = 0.00953477896100638

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.6710410000000001
$ws.Range("H16").Value = 2.013123
$ws.Range("I16").Value = 0.0572455703259029
$ws.Range("J16").Value = 0.08139652374383181
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 24.57775566666666
$ws.Range("N16").Value = 73.733267
$ws.Range("O16").Value = 0.1558145205313015
$ws.Range("P16").Value = 0.1827107372630203
$ws.Range("Q16").Value = 16.49268174031567
$ws.Range("R16").Value = 148.434135662841
$ws.Range("S16").Value = 0.00891969109287146
$ws.Range("T16").Value = 0.01487201886388245

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.6710410000000001
$ws.Range("H17").Value = 2.013123
$ws.Range("I17").Value = 0.0572455703259029
$ws.Range("J17").Value = 0.08139652374383181
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 14.32600733333333
$ws.Range("N17").Value = 42.978022
$ws.Range("O17").Value = 0.09082196088386706
$ws.Range("P17").Value = 0.1064993646046676
$ws.Range("Q17").Value = 9.613338286967334
$ws.Range("R17").Value = 86.520044582706
$ws.Range("S17").Value = 0.005199154948913814
$ws.Range("T17").Value = 0.008668678059746826

# Row 18
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 0.6710410000000001
$ws.Range("H18").Value = 2.013123
$ws.Range("I18").Value = 0.0572455703259029
$ws.Range("J18").Value = 0.08139652374383181
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 33.41628466666666
$ws.Range("N18").Value = 100.248854
$ws.Range("O18").Value = 0.2118477555025799
$ws.Range("P18").Value = 0.2484162545532246
$ws.Range("Q18").Value = 22.42369707900467
$ws.Range("R18").Value = 201.813273711042
$ws.Range("S18").Value = 0.01212734558600762
$ws.Range("T18").Value = 0.02022021956209532

# Row 19
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 0.6710410000000001
$ws.Range("H19").Value = 2.013123
$ws.Range("I19").Value = 0.0572455703259029
$ws.Range("J19").Value = 0.08139652374383181
$ws.Range("K19").Value = 2
$ws.Range("M19").Value = 5.4585215
$ws.Range("N19").Value = 10.917043
$ws.Range("O19").Value = 0.03460514954527787
$ws.Range("P19").Value = 0.02705238837799083
$ws.Range("Q19").Value = 3.6628917258815
$ws.Range("R19").Value = 21.977350355289
$ws.Range("S19").Value = 0.001980991521932591
$ws.Range("T19").Value = 0.00220197037293649

# Row 20
$ws.Range("E20").Value = 2
$ws.Range("G20").Value = 9.2083985
$ws.Range("H20").Value = 18.416797
$ws.Range("I20").Value = 0.7855556127281175
$ws.Range("J20").Value = 0.744645634815076
$ws.Range("K20").Value = 2
$ws.Range("M20").Value = 64.20135099999999
$ws.Range("N20").Value = 128.402702
$ws.Range("O20").Value = 0.4070144914449589
$ws.Range("P20").Value = 0.3181813759721767
$ws.Range("Q20").Value = 591.1916242463734
$ws.Range("R20").Value = 2364.766496985494
$ws.Range("S20").Value = 0.3197325182162679
$ws.Range("T20").Value = 0.2369323726971359

# Row 21
$ws.Range("E21").Value = 2
$ws.Range("G21").Value = 9.2083985
$ws.Range("H21").Value = 18.416797
$ws.Range("I21").Value = 0.7855556127281175
$ws.Range("J21").Value = 0.744645634815076
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 15.75734066666666
$ws.Range("N21").Value = 47.27202199999999
$ws.Range("O21").Value = 0.09989612209201491
$ws.Range("P21").Value = 0.11713987922892
$ws.Range("Q21").Value = 145.0998721589223
$ws.Range("R21").Value = 870.5992329535338
$ws.Range("S21").Value = 0.07847395939915561
$ws.Range("T21").Value = 0.08722769973058048

# Row 22
$ws.Range("E22").Value = 2
$ws.Range("G22").Value = 9.2083985
$ws.Range("H22").Value = 18.416797
$ws.Range("I22").Value = 0.7855556127281175
$ws.Range("J22").Value = 0.744645634815076
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 24.57775566666666
$ws.Range("N22").Value = 73.733267
$ws.Range("O22").Value = 0.1558145205313015
$ws.Range("P22").Value = 0.1827107372630203
$ws.Range("Q22").Value = 226.3217684142998
$ws.Range("R22").Value = 1357.930610485799
$ws.Range("S22").Value = 0.1224009711479044
$ws.Range("T22").Value = 0.1360547529367523

# Row 23
$ws.Range("E23").Value = 2
$ws.Range("G23").Value = 9.2083985
$ws.Range("H23").Value = 18.416797
$ws.Range("I23").Value = 0.7855556127281175
$ws.Range("J23").Value = 0.744645634815076
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 14.32600733333333
$ws.Range("N23").Value = 42.978022
$ws.Range("O23").Value = 0.09082196088386706
$ws.Range("P23").Value = 0.1064993646046676
$ws.Range("Q23").Value = 131.9195844392556
$ws.Range("R23").Value = 791.5175066355339
$ws.Range("S23").Value = 0.0713457011312953
$ws.Range("T23").Value = 0.07930428696344494

# Row 24
$ws.Range("E24").Value = 2
$ws.Range("G24").Value = 9.2083985
$ws.Range("H24").Value = 18.416797
$ws.Range("I24").Value = 0.7855556127281175
$ws.Range("J24").Value = 0.744645634815076
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 33.41628466666666
$ws.Range("N24").Value = 100.248854
$ws.Range("O24").Value = 0.2118477555025799
$ws.Range("P24").Value = 0.2484162545532246
$ws.Range("Q24").Value = 307.7104656001063
$ws.Range("R24").Value = 1846.262793600638
$ws.Range("S24").Value = 0.1664181933789056
$ws.Range("T24").Value = 0.1849820795701695

# Row 25
$ws.Range("E25").Value = 2
$ws.Range("G25").Value = 9.2083985
$ws.Range("H25").Value = 18.416797
$ws.Range("I25").Value = 0.7855556127281175
$ws.Range("J25").Value = 0.744645634815076
$ws.Range("K25").Value = 2
$ws.Range("M25").Value = 5.4585215
$ws.Range("N25").Value = 10.917043
$ws.Range("O25").Value = 0.03460514954527787
$ws.Range("P25").Value = 0.02705238837799083
$ws.Range("Q25").Value = 50.26424119281774
$ws.Range("R25").Value = 201.056964771271
$ws.Range("S25").Value = 0.02718426945458889
$ws.Range("T25").Value = 0.02014444291699296
